$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

$ws = $wb.Worksheets.Add($null, $ws1)
$ws.Name = "OpenAccountTest"

$ws.Range("A1").Value = "customer"
$ws.Range("B1").Value = "currency"
$ws.Range("A2").Value = "Huyen Ha"
$ws.Range("B2").Value = "VND"

$ws.Range("B2").Select() | Out-Null
